# Update the "想去人数" (column F) values across the four worksheets
# to match the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F
$updates = @{
    "展览"     = @{ 3=26574; 4=587; 5=253; 6=605; 8=544; 9=233; 10=357; 11=235; 14=299; 16=412; 18=1530; 19=197; 21=433; 22=100 }
    "演出"     = @{ 6=197; 7=197; 9=115; 10=115; 11=439; 16=58; 23=19 }
    "本地生活" = @{ 2=5045; 3=224 }
    "全部类型" = @{ 3=5045; 4=224; 5=26574; 6=587; 8=253; 10=605; 14=197; 15=197; 17=115; 18=115; 19=439; 20=544; 22=233; 23=357; 24=235; 28=299; 32=412; 34=58; 35=1530; 36=197; 39=433; 40=100; 47=19 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $value = $rowsForSheet[$row]
        $ws.Range("F$row").Value = $value
    }
}
